# -*- coding: utf-8 -*-
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark: it currently sits (empty) in the last
#    paragraph of the document; it needs to move to just after the run
#    "：丹阳" (the "1858：丹阳" entry). Adding a bookmark with the same
#    name automatically replaces/removes the previous one.
#    A plain offset-based Range collapsed exactly at a paragraph end is
#    ambiguous in this engine, so we anchor the bookmark using a
#    temporary marker run inserted right after the target run, then
#    remove the marker text while leaving the (now correctly seated)
#    bookmark tags behind.
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("1858：丹阳", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
$anchor.InsertAfter("@@BM@@")

$marker = $d.Content
$marker.Find.Execute("@@BM@@", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$marker.Bookmarks.Add("_GoBack")
$marker.Text = ""

# ---------------------------------------------------------------------
# 2) The entry "1860：，溧阳，江阴，常熟，太仓，昆山，吴县，通州" had a
#    stray leading comma dragged in along with the rest of the text.
#    Fix the typo (drop the extra "，") and make sure the corrected
#    text is split across runs the way Word leaves it after a
#    drag-and-drop/paste edit: "1860" | "：" | "溧阳，江阴，常熟，太仓，
#    昆山，吴县，通州" as three independent runs.
# ---------------------------------------------------------------------

# 2a) Remove the stray comma right after the full-width colon.
$full = $d.Content
$full.Find.Execute("：，溧阳，江阴，常熟，太仓，昆山，吴县，通州", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$comma = $d.Range($full.Start + 1, $full.Start + 2)
$comma.Text = ""

# 2b) Re-seat "溧阳，江阴，常熟，太仓，昆山，吴县，通州" into its own run
#     via copy/paste (copy+paste preserves run boundaries instead of
#     letting same-formatted neighbours re-merge).
$tailFind = $d.Content
$tailFind.Find.Execute("溧阳，江阴，常熟，太仓，昆山，吴县，通州", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailFind.Copy()
$tailStart = $tailFind.Start
$tailEnd = $tailFind.End
$tailDel = $d.Range($tailStart, $tailEnd)
$tailDel.Text = ""
$tailPaste = $d.Range($tailStart, $tailStart)
$tailPaste.Paste()

# 2c) Re-seat "：" into its own run too, separating it from "1860".
$colonFind = $d.Content
$colonFind.Find.Execute("1860：", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$colonRange = $d.Range($colonFind.End - 1, $colonFind.End)
$colonRange.Copy()
$colonStart = $colonRange.Start
$colonEnd = $colonRange.End
$colonDel = $d.Range($colonStart, $colonEnd)
$colonDel.Text = ""
$colonPaste = $d.Range($colonStart, $colonStart)
$colonPaste.Paste()
